# Auto-applies the numeric updates from the Odin_Profits.xlsx diff
# across all affected sheets/rows. Values were cross-checked against
# the original workbook cell-by-cell before generating this script.
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H9").Value = 1053.091
$ws.Range("I9").Value = 1246.7059
$ws.Range("J9").Value = 394.8
$ws.Range("K9").Value = 1246.7059
$ws.Range("L9").Value = 394.8
$ws.Range("M9").Value = -1077.7059
$ws.Range("N9").Value = -732.8
$ws.Range("H17").Value = 1334342.8
$ws.Range("J17").Value = 1334342.8
$ws.Range("L17").Value = 4003028.4
$ws.Range("N17").Value = -4003364.4
$ws.Range("H51").Value = 7196.643
$ws.Range("J51").Value = 7284.76
$ws.Range("L51").Value = 7284.76
$ws.Range("N51").Value = -8252.76
$ws.Range("H64").Value = 266669380
$ws.Range("I64").Value = 266669380
$ws.Range("K64").Value = 266669380
$ws.Range("M64").Value = -266669132
$ws.Range("H67").Value = 266669380
$ws.Range("I67").Value = 266669380
$ws.Range("K67").Value = 266669380
$ws.Range("M67").Value = -266668522
$ws.Range("H70").Value = 2365.8096
$ws.Range("J70").Value = 2360.2222
$ws.Range("L70").Value = 7080.6666
$ws.Range("N70").Value = -7620.6666
$ws.Range("H73").Value = 2365.8096
$ws.Range("J73").Value = 2360.2222
$ws.Range("L73").Value = 7080.6666
$ws.Range("N73").Value = -8952.6666
$ws.Range("H74").Value = 20431518
$ws.Range("I74").Value = 28581624
$ws.Range("J74").Value = 56250
$ws.Range("K74").Value = 28581624
$ws.Range("L74").Value = 56250
$ws.Range("M74").Value = -28580688
$ws.Range("N74").Value = -58122
$ws.Range("H77").Value = 20431518
$ws.Range("I77").Value = 28581624
$ws.Range("J77").Value = 56250
$ws.Range("K77").Value = 142908120
$ws.Range("L77").Value = 281250
$ws.Range("M77").Value = -142903440
$ws.Range("N77").Value = -290610
$ws.Range("H138").Value = 5330.3
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").ClearContents()

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 1373958.1
$ws.Range("I32").Value = 1021.2239
$ws.Range("J32").Value = 7944442
$ws.Range("K32").Value = 1021.2239
$ws.Range("L32").Value = 7944442
$ws.Range("M32").Value = -734.2239
$ws.Range("N32").Value = -7945016
$ws.Range("H63").Value = 2047.0834
$ws.Range("I63").Value = 1976.6
$ws.Range("J63").Value = 2399.5
$ws.Range("K63").Value = 1976.6
$ws.Range("L63").Value = 2399.5
$ws.Range("M63").Value = -1290.6
$ws.Range("N63").Value = -3771.5
$ws.Range("H66").Value = 2047.0834
$ws.Range("I66").Value = 1976.6
$ws.Range("J66").Value = 2399.5
$ws.Range("K66").Value = 9883
$ws.Range("L66").Value = 11997.5
$ws.Range("M66").Value = -6451
$ws.Range("N66").Value = -18861.5
$ws.Range("H92").Value = 1000000
$ws.Range("J92").Value = 1000000
$ws.Range("L92").Value = 1000000
$ws.Range("N92").Value = -1004992
$ws.Range("H94").Value = 800000
$ws.Range("J94").Value = 800000
$ws.Range("L94").Value = 800000
$ws.Range("N94").Value = -801802
$ws.Range("H97").Value = 1680.8572
$ws.Range("J97").Value = 2178.2
$ws.Range("L97").Value = 2178.2
$ws.Range("N97").Value = -3170.2
$ws.Range("H110").Value = 5432.64
$ws.Range("J110").Value = 7053.0586
$ws.Range("L110").Value = 7053.0586
$ws.Range("N110").Value = -11143.0586

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 5495610
$ws.Range("I20").Value = 6803842.5
$ws.Range("J20").Value = 1033
$ws.Range("K20").Value = 6803842.5
$ws.Range("L20").Value = 1033
$ws.Range("M20").Value = -6803595.5
$ws.Range("N20").Value = -1527
$ws.Range("H82").Value = 15569.429
$ws.Range("I82").Value = 14460.667
$ws.Range("J82").Value = 22222
$ws.Range("K82").Value = 14460.667
$ws.Range("L82").Value = 22222
$ws.Range("M82").Value = -14077.667
$ws.Range("N82").Value = -22988
$ws.Range("H85").Value = 15569.429
$ws.Range("I85").Value = 14460.667
$ws.Range("J85").Value = 22222
$ws.Range("K85").Value = 14460.667
$ws.Range("L85").Value = 22222
$ws.Range("M85").Value = -13134.667
$ws.Range("N85").Value = -24874
$ws.Range("H94").Value = 68722.12
$ws.Range("I94").Value = 5728.4
$ws.Range("K94").Value = 5728.4
$ws.Range("M94").Value = -5277.4
$ws.Range("H99").Value = 7495.394
$ws.Range("I99").Value = 6995.524
$ws.Range("K99").Value = 6995.524
$ws.Range("M99").Value = -5497.524

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 55559896
$ws.Range("J16").Value = 5374.25
$ws.Range("L16").Value = 5374.25
$ws.Range("N16").Value = -5948.25
$ws.Range("H31").Value = 4943.3076
$ws.Range("I31").Value = 2012.5555
$ws.Range("K31").Value = 2012.5555
$ws.Range("M31").Value = -1717.5555
$ws.Range("H34").Value = 4943.3076
$ws.Range("I34").Value = 2012.5555
$ws.Range("K34").Value = 2012.5555
$ws.Range("M34").Value = -1810.5555
$ws.Range("H92").Value = 90000
$ws.Range("I92").Value = 0
$ws.Range("K92").Value = 0
$ws.Range("M92").ClearContents()
$ws.Range("H94").Value = 43481732
$ws.Range("J94").Value = 5798.778
$ws.Range("L94").Value = 5798.778
$ws.Range("N94").Value = -6700.778
$ws.Range("H103").Value = 6333.3335
$ws.Range("I103").Value = 6333.3335
$ws.Range("K103").Value = 6333.3335
$ws.Range("M103").Value = -5161.3335
$ws.Range("H113").Value = 55559896
$ws.Range("J113").Value = 5374.25
$ws.Range("L113").Value = 5374.25
$ws.Range("N113").Value = -9714.25
$ws.Range("H133").Value = 33559.332
$ws.Range("I133").Value = 0
$ws.Range("K133").Value = 0
$ws.Range("M133").ClearContents()
$ws.Range("H134").Value = 50007620
$ws.Range("I134").Value = 62505476
$ws.Range("J134").Value = 16192.75
$ws.Range("K134").Value = 187516428
$ws.Range("L134").Value = 48578.25
$ws.Range("M134").Value = -187513893
$ws.Range("N134").Value = -53648.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H94").Value = 6949.524
$ws.Range("I94").Value = 6619.25
$ws.Range("J94").Value = 7027.2354
$ws.Range("K94").Value = 19857.75
$ws.Range("L94").Value = 21081.7062
$ws.Range("M94").Value = -19181.75
$ws.Range("N94").Value = -22433.7062
$ws.Range("H107").Value = 5022.857
$ws.Range("J107").Value = 5540.64
$ws.Range("L107").Value = 16621.92
$ws.Range("N107").Value = -20461.92
$ws.Range("H134").Value = 21706.785
$ws.Range("I134").Value = 18626.818
$ws.Range("K134").Value = 55880.454
$ws.Range("M134").Value = -50810.454
$ws.Range("H141").Value = 2747.923
$ws.Range("I141").Value = 2747.923
$ws.Range("K141").Value = 8243.769
$ws.Range("M141").Value = -3063.769

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H47").Value = 30000
$ws.Range("J47").Value = 30000
$ws.Range("L47").Value = 30000
$ws.Range("N47").Value = -31136
$ws.Range("H49").Value = 6500
$ws.Range("J49").Value = 6500
$ws.Range("L49").Value = 6500
$ws.Range("N49").Value = -6868
$ws.Range("H113").Value = 8079.125
$ws.Range("J113").Value = 9345.846
$ws.Range("L113").Value = 9345.846
$ws.Range("N113").Value = -13685.846
$ws.Range("H122").Value = 3842.372
$ws.Range("I122").Value = 2970.5806
$ws.Range("J122").Value = 6094.5
$ws.Range("K122").Value = 8911.7418
$ws.Range("L122").Value = 18283.5
$ws.Range("M122").Value = -6461.7418
$ws.Range("N122").Value = -23183.5

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 38462530
$ws.Range("I46").Value = 1063.3334
$ws.Range("J46").Value = 71429496
$ws.Range("K46").Value = 1063.3334
$ws.Range("L46").Value = 71429496
$ws.Range("M46").Value = -875.3334
$ws.Range("N46").Value = -71429872

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 25633.334
$ws.Range("I62").Value = 17075
$ws.Range("J62").Value = 42750
$ws.Range("K62").Value = 17075
$ws.Range("L62").Value = 42750
$ws.Range("M62").Value = -16451
$ws.Range("N62").Value = -43998
$ws.Range("H65").Value = 25633.334
$ws.Range("I65").Value = 17075
$ws.Range("J65").Value = 42750
$ws.Range("K65").Value = 85375
$ws.Range("L65").Value = 213750
$ws.Range("M65").Value = -82255
$ws.Range("N65").Value = -219990
$ws.Range("H81").Value = 1280.6897
$ws.Range("I81").Value = 1226.909
$ws.Range("J81").Value = 1449.7142
$ws.Range("K81").Value = 2453.818
$ws.Range("L81").Value = 2899.4284
$ws.Range("M81").Value = -1392.818
$ws.Range("N81").Value = -5021.4284
$ws.Range("H84").Value = 1280.6897
$ws.Range("I84").Value = 1226.909
$ws.Range("J84").Value = 1449.7142
$ws.Range("K84").Value = 12269.09
$ws.Range("L84").Value = 14497.142
$ws.Range("M84").Value = -6965.09
$ws.Range("N84").Value = -25105.142
$ws.Range("H126").Value = 3219.0286
$ws.Range("I126").Value = 2404.138
$ws.Range("K126").Value = 7212.414
$ws.Range("M126").Value = -4742.414
$ws.Range("H132").Value = 7130.3794
$ws.Range("I132").Value = 3872.72
$ws.Range("K132").Value = 11618.16
$ws.Range("M132").Value = -9088.16
